$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C26").Value = "Lecture 15: Confidence Regions and Intervals"
$ws.Range("C20").Value = "Lecture 13: Testing Hypotheses"
$ws.Range("C21").Value = "Lecture 14: Testing Submodels"
$ws.Range("C27").Value = "Lecture 16:"
$ws.Range("C29").Value = "Lecture 17:"
$ws.Range("C30").Value = "Lecture 18:"
$ws.Range("C32").Value = "Lecture 19: "
$ws.Range("C33").Value = "Lecture 20:"
$ws.Range("D26").Value = "15-CI"

$ws.Range("D26").Select()
